$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "SD008-Ajicero"
$ws.Range("B2").Value = "SD008"
$ws.Range("G2").Value = 2
$ws.Range("P2").Value = "SD008-Ajicero"
$ws.Range("T2").Value = "SD008"
